# Kazakhstan Premier League workbook update
# - Swap a handful of existing match rows (B:AC) back into their correct
#   chronological positions (ids in column A stay fixed).
# - Append four newly scraped matches as rows 111-114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row swaps: exchange everything except column A (the running index)
#    between each pair of rows.
# ---------------------------------------------------------------------
$swapPairs = @(
    @(9, 10),
    @(63, 64),
    @(85, 86),
    @(92, 93),
    @(98, 100)
)

foreach ($pair in $swapPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

# ---------------------------------------------------------------------
# 2) Append four new match rows (111-114), copying the formatting from
#    the last existing data row (110) so styles/number-formats match.
# ---------------------------------------------------------------------
$lastRow = 110
$newRows = @(
    @{ Row=111; A=109; B=7914794; E=45356.375;            F="Zhetysu";          G="FC Astana";        H=0; I=2; J="A";
       K=6.5;   L=4;    M=1.4;    N=7;    O=3.75; P=1.444; Q=1.25;  R=1.8;   S=2;     T=2.25;
       U=1.8;   V=2;    W=-1;     X=-1;   Y=0.444; Z=-1;    AA=1;    AB=-0.5; AC=0.5 },
    @{ Row=112; A=110; B=7874788; E=45357.27083333334;    F="Kaisar Kyzylorda"; G="FK Kyzylzhar";      H=0; I=3; J="A";
       K=2.375; L=3.2;  M=2.625;  N=2.55; O=3;    P=2.55;  Q=0;     R=1.875; S=1.925; T=1.75;
       U=1.825; V=1.975;W=-1;     X=-1;   Y=1.55;  Z=-1;    AA=0.925;AB=0.825; AC=-1 },
    @{ Row=113; A=111; B=7874789; E=45357.375;            F="FK Atyrau";        G="Ordabasy";         H=1; I=1; J="D";
       K=4.333; L=3.25; M=1.727;  N=5.75; O=3.3;  P=1.55;  Q=1;     R=1.8;   S=2;     T=2.25;
       U=2.025; V=1.775;W=-1;     X=2.3;  Y=-1;    Z=0.8;   AA=-1;   AB=-0.5; AC=0.3875 },
    @{ Row=114; A=112; B=7874790; E=45357.47916666666;    F="Tobol Kostanay";   G="FC Elimai Semey";  H=4; I=2; J="H";
       K=1.666; L=3.4;  M=4.5;    N=1.5;  O=3.6;  P=5.75;  Q=-1;    R=1.9;   S=1.9;   T=2.25;
       U=1.825; V=1.975;W=0.5;    X=-1;   Y=-1;    Z=0.8999999999999999; AA=-1; AB=0.825; AC=-1 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Clone formatting (number formats, borders, bold id column, etc.)
    # from the previous row before touching any values.
    $ws.Range("A$lastRow`:AC$lastRow").Copy() | Out-Null
    $ws.Range("A$r`:AC$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value2 = $nr.A
    $ws.Cells.Item($r, 2).Value2 = $nr.B
    $ws.Cells.Item($r, 3).Value2 = "Kazakhstan Premier League"
    $ws.Cells.Item($r, 4).Value2 = "Kazakhstan Premier League"
    $ws.Cells.Item($r, 5).Value2 = $nr.E
    $ws.Cells.Item($r, 6).Value2 = $nr.F
    $ws.Cells.Item($r, 7).Value2 = $nr.G
    $ws.Cells.Item($r, 8).Value2 = $nr.H
    $ws.Cells.Item($r, 9).Value2 = $nr.I
    $ws.Cells.Item($r, 10).Value2 = $nr.J
    $ws.Cells.Item($r, 11).Value2 = $nr.K
    $ws.Cells.Item($r, 12).Value2 = $nr.L
    $ws.Cells.Item($r, 13).Value2 = $nr.M
    $ws.Cells.Item($r, 14).Value2 = $nr.N
    $ws.Cells.Item($r, 15).Value2 = $nr.O
    $ws.Cells.Item($r, 16).Value2 = $nr.P
    $ws.Cells.Item($r, 17).Value2 = $nr.Q
    $ws.Cells.Item($r, 18).Value2 = $nr.R
    $ws.Cells.Item($r, 19).Value2 = $nr.S
    $ws.Cells.Item($r, 20).Value2 = $nr.T
    $ws.Cells.Item($r, 21).Value2 = $nr.U
    $ws.Cells.Item($r, 22).Value2 = $nr.V
    $ws.Cells.Item($r, 23).Value2 = $nr.W
    $ws.Cells.Item($r, 24).Value2 = $nr.X
    $ws.Cells.Item($r, 25).Value2 = $nr.Y
    $ws.Cells.Item($r, 26).Value2 = $nr.Z
    $ws.Cells.Item($r, 27).Value2 = $nr.AA
    $ws.Cells.Item($r, 28).Value2 = $nr.AB
    $ws.Cells.Item($r, 29).Value2 = $nr.AC

    $lastRow = $r
}
